$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.471.77'
$ws.Range('E2').Value = '  +9.19%  '
$ws.Range('D3').Value = '1.681.80'
$ws.Range('E3').Value = '  +4.78%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9989'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3708'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3443'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.66'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +12.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.178'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07275'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.162'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.38'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.756'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.68%  '
$ws.Range('D16').Value = '1.679.28'
$ws.Range('E16').Value = '  +4.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001109'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9989'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06696'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '81.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.119'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.26'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.47%  '
$ws.Range('D24').Value = '24.428.59'
$ws.Range('E24').Value = '  +8.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.443'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.665'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '152.69'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.49'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('D29').Value = '1.865.03'
$ws.Range('E29').Value = '  +4.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.356'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.050'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9776'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.91%  '
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08472'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.19%  '
$ws.Range('B35').Value = 'WEMIXTOKEN'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.709'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.47'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06519'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.386'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.904'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02336'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.91%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.264'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2116'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6179'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.38%  '
$ws.Range('E44').Value = '  +0.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.20'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.49%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.774'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.14%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5965'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '127.38'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.029'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07199'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.81%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '76.08'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.08%  '
